# Insert two new price records (rows) right before the current row 343,
# pushing the existing rows 343-360 down to 345-362.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 343 (old row 343 becomes 345, etc.)
$ws.Rows.Item(343).Resize(2).Insert()

# --- New row 343 ---
$ws.Cells.Item(343, 1).Value = 10
$ws.Cells.Item(343, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(343, 3).Value = "La Araucanía"
$ws.Cells.Item(343, 4).Value = 45041
$ws.Cells.Item(343, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(343, 5).Value = 9
$ws.Cells.Item(343, 6).Value = "Fruta"
$ws.Cells.Item(343, 7).Value = 100103
$ws.Cells.Item(343, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(343, 9).Value = 100103002
$ws.Cells.Item(343, 10).Value = "Ciruela"
$ws.Cells.Item(343, 11).Value = "Blue Giant"
$ws.Cells.Item(343, 12).Value = "Primera"
$ws.Cells.Item(343, 13).Value = 125
$ws.Cells.Item(343, 14).Value = 14000
$ws.Cells.Item(343, 15).Value = 14000
$ws.Cells.Item(343, 16).Value = 14000
$ws.Cells.Item(343, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(343, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(343, 19).Value = 778
$ws.Cells.Item(343, 20).Value = 18

# --- New row 344 ---
$ws.Cells.Item(344, 1).Value = 10
$ws.Cells.Item(344, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(344, 3).Value = "La Araucanía"
$ws.Cells.Item(344, 4).Value = 45041
$ws.Cells.Item(344, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(344, 5).Value = 9
$ws.Cells.Item(344, 6).Value = "Fruta"
$ws.Cells.Item(344, 7).Value = 100103
$ws.Cells.Item(344, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(344, 9).Value = 100103002
$ws.Cells.Item(344, 10).Value = "Ciruela"
$ws.Cells.Item(344, 11).Value = "Blue Giant"
$ws.Cells.Item(344, 12).Value = "Segunda"
$ws.Cells.Item(344, 13).Value = 65
$ws.Cells.Item(344, 14).Value = 10000
$ws.Cells.Item(344, 15).Value = 10000
$ws.Cells.Item(344, 16).Value = 10000
$ws.Cells.Item(344, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(344, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(344, 19).Value = 556
$ws.Cells.Item(344, 20).Value = 18
